$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.916.96"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "1.910.04"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5047"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4050"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08270"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.100"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.77%  "
$ws.Range("D13").Value = "1.907.11"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.385"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.209"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06506"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.72%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "29.960.96"
$ws.Range("E23").Value = "  -0.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.196"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Value = "2.128.17"
$ws.Range("E27").Value = "  +0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.294"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.124"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("E32").Value = "  -2.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.923"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.805"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.412"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02436"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06371"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2145"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.195"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.691"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6457"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.211"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.214"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6032"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.633"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "78.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.125"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.36%  "
